$wb = $excel.ActiveWorkbook

# Rename sheet "Groeperen" to "Groepen"
$groepenSheet = $wb.Worksheets.Item("Groeperen")
$groepenSheet.Name = "Groepen"

# Update the selection on the "Groepen" sheet (was B10, now B17), and make it the active/selected sheet
$groepenSheet.Activate()
$groepenSheet.Range("B17").Select()

# The "Irrelevant" sheet keeps its own selection at A3 (no change needed there)
